# Apply changes described by the diff across 4 worksheets:
# Matriz_Resultados, P_valores, Estadisticos_HLN_DM, Resumen_Modelos
$wb = $excel.ActiveWorkbook

# --- Sheet: Matriz_Resultados ---
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("D3").Value = 0
$wsMatriz.Range("C4").Value = 0

# --- Sheet: P_valores (plain numeric cells) ---
$wsP = $wb.Worksheets.Item("P_valores")
$wsP.Range("B2").Value = 1
$wsP.Range("C2").Value = [double]"7.620427444621214E-08"
$wsP.Range("D2").Value = [double]"7.637692323037015E-08"
$wsP.Range("E2").Value = 0.005305285977035057
$wsP.Range("F2").Value = 0.004131548874835911
$wsP.Range("G2").Value = 0.05111695843634334
$wsP.Range("H2").Value = 0.00412467555736673
$wsP.Range("I2").Value = 0.005424700392069992
$wsP.Range("J2").Value = 0.005647142000966943
$wsP.Range("B3").Value = [double]"7.620427444621214E-08"
$wsP.Range("C3").Value = 1
$wsP.Range("D3").Value = 0.1674540400586997
$wsP.Range("E3").Value = 0.02238178758511578
$wsP.Range("F3").Value = 0.01743993615042294
$wsP.Range("G3").Value = 0.1051148974766396
$wsP.Range("H3").Value = 0.008548870124466923
$wsP.Range("I3").Value = 0.008166489730020876
$wsP.Range("J3").Value = 0.008090275648662271
$wsP.Range("B4").Value = [double]"7.637692323037015E-08"
$wsP.Range("C4").Value = 0.1674540400586997
$wsP.Range("D4").Value = 1
$wsP.Range("E4").Value = 0.02393143720546975
$wsP.Range("F4").Value = 0.01871609137959562
$wsP.Range("G4").Value = 0.1080827303257397
$wsP.Range("H4").Value = 0.008875912317723778
$wsP.Range("I4").Value = 0.008330163298001292
$wsP.Range("J4").Value = 0.008233965610212612
$wsP.Range("B5").Value = 0.005305285977035057
$wsP.Range("C5").Value = 0.02238178758511578
$wsP.Range("D5").Value = 0.02393143720546975
$wsP.Range("E5").Value = 1
$wsP.Range("F5").Value = 0.4863314953090097
$wsP.Range("G5").Value = 0.3840565803825529
$wsP.Range("H5").Value = 0.006038550104455975
$wsP.Range("I5").Value = 0.005629114665586243
$wsP.Range("J5").Value = 0.005904426321486245
$wsP.Range("B6").Value = 0.004131548874835911
$wsP.Range("C6").Value = 0.01743993615042294
$wsP.Range("D6").Value = 0.01871609137959562
$wsP.Range("E6").Value = 0.4863314953090097
$wsP.Range("F6").Value = 1
$wsP.Range("G6").Value = 0.4488142296517927
$wsP.Range("H6").Value = 0.007730055059711249
$wsP.Range("I6").Value = 0.006341787075952121
$wsP.Range("J6").Value = 0.006524913545280064
$wsP.Range("B7").Value = 0.05111695843634334
$wsP.Range("C7").Value = 0.1051148974766396
$wsP.Range("D7").Value = 0.1080827303257397
$wsP.Range("E7").Value = 0.3840565803825529
$wsP.Range("F7").Value = 0.4488142296517927
$wsP.Range("G7").Value = 1
$wsP.Range("H7").Value = 0.02546020839792718
$wsP.Range("I7").Value = 0.0005354153139132123
$wsP.Range("J7").Value = 0.0009665947134080977
$wsP.Range("B8").Value = 0.00412467555736673
$wsP.Range("C8").Value = 0.008548870124466923
$wsP.Range("D8").Value = 0.008875912317723778
$wsP.Range("E8").Value = 0.006038550104455975
$wsP.Range("F8").Value = 0.007730055059711249
$wsP.Range("G8").Value = 0.02546020839792718
$wsP.Range("H8").Value = 1
$wsP.Range("I8").Value = 0.01493134748401359
$wsP.Range("J8").Value = 0.01086650606326289
$wsP.Range("B9").Value = 0.005424700392069992
$wsP.Range("C9").Value = 0.008166489730020876
$wsP.Range("D9").Value = 0.008330163298001292
$wsP.Range("E9").Value = 0.005629114665586243
$wsP.Range("F9").Value = 0.006341787075952121
$wsP.Range("G9").Value = 0.0005354153139132123
$wsP.Range("H9").Value = 0.01493134748401359
$wsP.Range("I9").Value = 1
$wsP.Range("J9").Value = 0.01376997561758664
$wsP.Range("B10").Value = 0.005647142000966943
$wsP.Range("C10").Value = 0.008090275648662271
$wsP.Range("D10").Value = 0.008233965610212612
$wsP.Range("E10").Value = 0.005904426321486245
$wsP.Range("F10").Value = 0.006524913545280064
$wsP.Range("G10").Value = 0.0009665947134080977
$wsP.Range("H10").Value = 0.01086650606326289
$wsP.Range("I10").Value = 0.01376997561758664
$wsP.Range("J10").Value = 1

# --- Sheet: Estadisticos_HLN_DM (plain numeric cells) ---
$wsE = $wb.Worksheets.Item("Estadisticos_HLN_DM")
$wsE.Range("B2").Value = 0
$wsE.Range("C2").Value = -6.817818191409489
$wsE.Range("D2").Value = -6.817047592620416
$wsE.Range("E2").Value = -2.979011395689851
$wsE.Range("F2").Value = -3.075301860208742
$wsE.Range("G2").Value = -2.021860288360768
$wsE.Range("H2").Value = -3.075938670891716
$wsE.Range("I2").Value = -2.970374700405864
$wsE.Range("J2").Value = -2.954753774119182
$wsE.Range("B3").Value = 6.817818191409489
$wsE.Range("C3").Value = 0
$wsE.Range("D3").Value = -1.410575500942601
$wsE.Range("E3").Value = -2.392943430196895
$wsE.Range("F3").Value = -2.499190314855062
$wsE.Range("G3").Value = -1.664939680568258
$wsE.Range("H3").Value = -2.791336674250943
$wsE.Range("I3").Value = -2.809582320456876
$wsE.Range("J3").Value = -2.813314137882683
$wsE.Range("B4").Value = 6.817047592620416
$wsE.Range("C4").Value = 1.410575500942601
$wsE.Range("D4").Value = 0
$wsE.Range("E4").Value = -2.364024560027169
$wsE.Range("F4").Value = -2.469350628913434
$wsE.Range("G4").Value = -1.65033987872114
$wsE.Range("H4").Value = -2.776326266125144
$wsE.Range("I4").Value = -2.801676823439771
$wsE.Range("J4").Value = -2.806305432527586
$wsE.Range("B5").Value = 2.979011395689851
$wsE.Range("C5").Value = 2.392943430196895
$wsE.Range("D5").Value = 2.364024560027169
$wsE.Range("E5").Value = 0
$wsE.Range("F5").Value = -0.7038296828992813
$wsE.Range("G5").Value = -0.8818376846265282
$wsE.Range("H5").Value = -2.928623257780381
$wsE.Range("I5").Value = -2.955997956507879
$wsE.Range("J5").Value = -2.937393079054316
$wsE.Range("B6").Value = 3.075301860208742
$wsE.Range("C6").Value = 2.499190314855062
$wsE.Range("D6").Value = 2.469350628913434
$wsE.Range("E6").Value = 0.7038296828992813
$wsE.Range("F6").Value = 0
$wsE.Range("G6").Value = -0.7662515165628037
$wsE.Range("H6").Value = -2.831409553940748
$wsE.Range("I6").Value = -2.90945244358324
$wsE.Range("J6").Value = -2.898288175503464
$wsE.Range("B7").Value = 2.021860288360768
$wsE.Range("C7").Value = 1.664939680568258
$wsE.Range("D7").Value = 1.65033987872114
$wsE.Range("E7").Value = 0.8818376846265282
$wsE.Range("F7").Value = 0.7662515165628037
$wsE.Range("G7").Value = 0
$wsE.Range("H7").Value = -2.337111901495963
$wsE.Range("I7").Value = -3.823486950221935
$wsE.Range("J7").Value = -3.612944267251107
$wsE.Range("B8").Value = 3.075938670891716
$wsE.Range("C8").Value = 2.791336674250943
$wsE.Range("D8").Value = 2.776326266125144
$wsE.Range("E8").Value = 2.928623257780381
$wsE.Range("F8").Value = 2.831409553940748
$wsE.Range("G8").Value = 2.337111901495963
$wsE.Range("H8").Value = 0
$wsE.Range("I8").Value = -2.564193228957207
$wsE.Range("J8").Value = -2.694754620437258
$wsE.Range("B9").Value = 2.970374700405864
$wsE.Range("C9").Value = 2.809582320456876
$wsE.Range("D9").Value = 2.801676823439771
$wsE.Range("E9").Value = 2.955997956507879
$wsE.Range("F9").Value = 2.90945244358324
$wsE.Range("G9").Value = 3.823486950221935
$wsE.Range("H9").Value = 2.564193228957207
$wsE.Range("I9").Value = 0
$wsE.Range("J9").Value = -2.597763446613827
$wsE.Range("B10").Value = 2.954753774119182
$wsE.Range("C10").Value = 2.813314137882683
$wsE.Range("D10").Value = 2.806305432527586
$wsE.Range("E10").Value = 2.937393079054316
$wsE.Range("F10").Value = 2.898288175503464
$wsE.Range("G10").Value = 3.612944267251107
$wsE.Range("H10").Value = 2.694754620437258
$wsE.Range("I10").Value = 2.597763446613827
$wsE.Range("J10").Value = 0

# --- Sheet: Resumen_Modelos ---
$wsR = $wb.Worksheets.Item("Resumen_Modelos")

# Helper: write a numeric-looking / percent string as TEXT (preserve formatting
# like trailing zeros / % sign) without leaving a visible style on the cell.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$wsR.Range("A2").Value = "Sieve Bootstrap"
$wsR.Range("B2").Value = 2
$wsR.Range("C2").Value = 0
$wsR.Range("D2").Value = 6
Set-TextValue $wsR.Range("E2") "25.0%"
Set-TextValue $wsR.Range("F2") "0.5705"
Set-TextValue $wsR.Range("G2") "0.3071"
Set-TextValue $wsR.Range("H2") "0.5383"

$wsR.Range("A3").Value = "LSPMW"
$wsR.Range("B3").Value = 0
$wsR.Range("C3").Value = 1
$wsR.Range("D3").Value = 7
Set-TextValue $wsR.Range("E3") "0.0%"
Set-TextValue $wsR.Range("F3") "0.8362"
Set-TextValue $wsR.Range("G3") "0.6845"
Set-TextValue $wsR.Range("H3") "0.8186"

$wsR.Range("A4").Value = "LSPM"
$wsR.Range("B4").Value = 0
$wsR.Range("C4").Value = 1
$wsR.Range("D4").Value = 7
Set-TextValue $wsR.Range("E4") "0.0%"
Set-TextValue $wsR.Range("F4") "0.8455"
Set-TextValue $wsR.Range("G4") "0.7015"
Set-TextValue $wsR.Range("H4") "0.8297"

$wsR.Range("A5").Value = "MCPS"
$wsR.Range("B5").Value = 0
$wsR.Range("C5").Value = 0
$wsR.Range("D5").Value = 8
Set-TextValue $wsR.Range("E5") "0.0%"
Set-TextValue $wsR.Range("F5") "1.5627"
Set-TextValue $wsR.Range("G5") "3.5113"
Set-TextValue $wsR.Range("H5") "2.2469"

$wsR.Range("A6").Value = "AV-MCPS"
$wsR.Range("B6").Value = 0
$wsR.Range("C6").Value = 0
$wsR.Range("D6").Value = 8
Set-TextValue $wsR.Range("E6") "0.0%"
Set-TextValue $wsR.Range("F6") "1.5887"
Set-TextValue $wsR.Range("G6") "3.3264"
Set-TextValue $wsR.Range("H6") "2.0937"

$wsR.Range("A7").Value = "DeepAR"
$wsR.Range("B7").Value = 2
$wsR.Range("C7").Value = 0
$wsR.Range("D7").Value = 6
Set-TextValue $wsR.Range("E7") "25.0%"
Set-TextValue $wsR.Range("F7") "1.8360"
Set-TextValue $wsR.Range("G7") "7.4159"
Set-TextValue $wsR.Range("H7") "4.0391"

$wsR.Range("A8").Value = "EnCQR-LSTM"
$wsR.Range("B8").Value = 0
$wsR.Range("C8").Value = 0
$wsR.Range("D8").Value = 8
Set-TextValue $wsR.Range("E8") "0.0%"
Set-TextValue $wsR.Range("F8") "2.4549"
Set-TextValue $wsR.Range("G8") "5.4040"
Set-TextValue $wsR.Range("H8") "2.2013"

$wsR.Range("A9").Value = "AREPD"
$wsR.Range("B9").Value = 0
$wsR.Range("C9").Value = 1
$wsR.Range("D9").Value = 7
Set-TextValue $wsR.Range("E9") "0.0%"
Set-TextValue $wsR.Range("F9") "3.8803"
Set-TextValue $wsR.Range("G9") "11.4195"
Set-TextValue $wsR.Range("H9") "2.9430"

$wsR.Range("A10").Value = "Block Bootstrapping"
$wsR.Range("B10").Value = 0
$wsR.Range("C10").Value = 1
$wsR.Range("D10").Value = 7
Set-TextValue $wsR.Range("E10") "0.0%"
Set-TextValue $wsR.Range("F10") "4.2755"
Set-TextValue $wsR.Range("G10") "12.2657"
Set-TextValue $wsR.Range("H10") "2.8688"

